$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 96
$ws.Cells.Item(96, 8).Value = 1391.4667
$ws.Cells.Item(96, 9).Value = 958.5
$ws.Cells.Item(96, 10).Value = 1680.1111
$ws.Cells.Item(96, 11).Value = 2875.5
$ws.Cells.Item(96, 12).Value = 5040.3333
$ws.Cells.Item(96, 13).Value = -1502.5
$ws.Cells.Item(96, 14).Value = -7786.3333
# Row 101
$ws.Cells.Item(101, 8).Value = 293.33334
$ws.Cells.Item(101, 9).Value = 180
$ws.Cells.Item(101, 10).Value = 350
$ws.Cells.Item(101, 11).Value = 540
$ws.Cells.Item(101, 12).Value = 1050
$ws.Cells.Item(101, 13).Value = 1082
$ws.Cells.Item(101, 14).Value = -4294
# Row 116
$ws.Cells.Item(116, 8).Value = 12503990
$ws.Cells.Item(116, 9).Value = 22728700
$ws.Cells.Item(116, 11).Value = 22728700
$ws.Cells.Item(116, 13).Value = -22725258
# Row 129
$ws.Cells.Item(129, 8).Value = 753.3019
$ws.Cells.Item(129, 10).Value = 800.8542
$ws.Cells.Item(129, 12).Value = 2402.5626
$ws.Cells.Item(129, 14).Value = -12402.5626
# Row 132
$ws.Cells.Item(132, 8).Value = 4383.913
$ws.Cells.Item(132, 9).Value = 4852.778
$ws.Cells.Item(132, 10).Value = 2696
$ws.Cells.Item(132, 11).Value = 14558.334
$ws.Cells.Item(132, 12).Value = 8088
$ws.Cells.Item(132, 13).Value = -12028.334
$ws.Cells.Item(132, 14).Value = -13148

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Cells.Item(32, 8).Value = 6223.1445
$ws.Cells.Item(32, 9).Value = 4752.683
$ws.Cells.Item(32, 10).Value = 14261.667
$ws.Cells.Item(32, 11).Value = 4752.683
$ws.Cells.Item(32, 12).Value = 14261.667
$ws.Cells.Item(32, 13).Value = -4465.683
$ws.Cells.Item(32, 14).Value = -14835.667
# Row 61
$ws.Cells.Item(61, 8).Value = 2074.5454
$ws.Cells.Item(61, 9).Value = 1729.3077
$ws.Cells.Item(61, 10).Value = 3356.8572
$ws.Cells.Item(61, 11).Value = 1729.3077
$ws.Cells.Item(61, 12).Value = 3356.8572
$ws.Cells.Item(61, 13).Value = -1517.3077
$ws.Cells.Item(61, 14).Value = -3780.8572
# Row 74
$ws.Cells.Item(74, 8).Value = 32259668
$ws.Cells.Item(74, 9).Value = 45455230
$ws.Cells.Item(74, 10).Value = 3840
$ws.Cells.Item(74, 11).Value = 45455230
$ws.Cells.Item(74, 12).Value = 3840
$ws.Cells.Item(74, 13).Value = -45454356
$ws.Cells.Item(74, 14).Value = -5588
# Row 77
$ws.Cells.Item(77, 8).Value = 32259668
$ws.Cells.Item(77, 9).Value = 45455230
$ws.Cells.Item(77, 10).Value = 3840
$ws.Cells.Item(77, 11).Value = 227276150
$ws.Cells.Item(77, 12).Value = 19200
$ws.Cells.Item(77, 13).Value = -227271782
$ws.Cells.Item(77, 14).Value = -27936
# Row 132
$ws.Cells.Item(132, 8).Value = 15967.081
$ws.Cells.Item(132, 9).Value = 2285.1333
$ws.Cells.Item(132, 10).Value = 74604
$ws.Cells.Item(132, 11).Value = 6855.3999
$ws.Cells.Item(132, 12).Value = 223812
$ws.Cells.Item(132, 13).Value = -4325.3999
$ws.Cells.Item(132, 14).Value = -228872
# Row 136
$ws.Cells.Item(136, 8).Value = 2074.5454
$ws.Cells.Item(136, 9).Value = 1729.3077
$ws.Cells.Item(136, 10).Value = 3356.8572
$ws.Cells.Item(136, 11).Value = 5187.9231
$ws.Cells.Item(136, 12).Value = 10070.5716
$ws.Cells.Item(136, 13).Value = -2637.9231
$ws.Cells.Item(136, 14).Value = -15170.5716

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Cells.Item(20, 8).Value = 1839.8695
$ws.Cells.Item(20, 9).Value = 1864.2
$ws.Cells.Item(20, 10).Value = 1821.1538
$ws.Cells.Item(20, 11).Value = 1864.2
$ws.Cells.Item(20, 12).Value = 1821.1538
$ws.Cells.Item(20, 13).Value = -1617.2
$ws.Cells.Item(20, 14).Value = -2315.1538
# Row 25
$ws.Cells.Item(25, 8).Value = 484
$ws.Cells.Item(25, 9).Value = 484
$ws.Cells.Item(25, 10).Value = 0
$ws.Cells.Item(25, 11).Value = 484
$ws.Cells.Item(25, 12).Value = 0
$ws.Cells.Item(25, 13).Value = -249
$ws.Cells.Item(25, 14).ClearContents()
# Row 54
$ws.Cells.Item(54, 8).Value = 8193.25
$ws.Cells.Item(54, 9).Value = 2909.2
$ws.Cells.Item(54, 11).Value = 2909.2
$ws.Cells.Item(54, 13).Value = -2425.2
# Row 134
$ws.Cells.Item(134, 8).Value = 3739.1143
$ws.Cells.Item(134, 9).Value = 3886.5151
$ws.Cells.Item(134, 11).Value = 11659.5453
$ws.Cells.Item(134, 13).Value = -9124.5453

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Cells.Item(16, 8).Value = 1172.4166
$ws.Cells.Item(16, 9).Value = 911.6667
$ws.Cells.Item(16, 10).Value = 1433.1666
$ws.Cells.Item(16, 11).Value = 911.6667
$ws.Cells.Item(16, 12).Value = 1433.1666
$ws.Cells.Item(16, 13).Value = -624.6667
$ws.Cells.Item(16, 14).Value = -2007.1666
# Row 22
$ws.Cells.Item(22, 8).Value = 174.875
$ws.Cells.Item(22, 9).Value = 153.6923
$ws.Cells.Item(22, 10).Value = 266.66666
$ws.Cells.Item(22, 11).Value = 153.6923
$ws.Cells.Item(22, 12).Value = 266.66666
$ws.Cells.Item(22, 13).Value = 196.3077
$ws.Cells.Item(22, 14).Value = -966.66666
# Row 58
$ws.Cells.Item(58, 8).Value = 17045.344
$ws.Cells.Item(58, 9).Value = 1882.7693
$ws.Cells.Item(58, 10).Value = 27419.736
$ws.Cells.Item(58, 11).Value = 1882.7693
$ws.Cells.Item(58, 12).Value = 27419.736
$ws.Cells.Item(58, 13).Value = -1679.7693
$ws.Cells.Item(58, 14).Value = -27825.736
# Row 94
$ws.Cells.Item(94, 8).Value = 3585.125
$ws.Cells.Item(94, 10).Value = 4217.5557
$ws.Cells.Item(94, 12).Value = 4217.5557
$ws.Cells.Item(94, 14).Value = -5119.5557
# Row 105
$ws.Cells.Item(105, 8).Value = 746.9167
$ws.Cells.Item(105, 9).Value = 746.9167
$ws.Cells.Item(105, 11).Value = 746.9167
$ws.Cells.Item(105, 13).Value = 1000.0833
# Row 113
$ws.Cells.Item(113, 8).Value = 1172.4166
$ws.Cells.Item(113, 9).Value = 911.6667
$ws.Cells.Item(113, 10).Value = 1433.1666
$ws.Cells.Item(113, 11).Value = 911.6667
$ws.Cells.Item(113, 12).Value = 1433.1666
$ws.Cells.Item(113, 13).Value = 1258.3333
$ws.Cells.Item(113, 14).Value = -5773.1666
# Row 134
$ws.Cells.Item(134, 8).Value = 947.8
$ws.Cells.Item(134, 9).Value = 889.2105
$ws.Cells.Item(134, 10).Value = 1133.3334
$ws.Cells.Item(134, 11).Value = 2667.6315
$ws.Cells.Item(134, 12).Value = 3400.0002
$ws.Cells.Item(134, 13).Value = -132.6315
$ws.Cells.Item(134, 14).Value = -8470.0002
# Row 136
$ws.Cells.Item(136, 8).Value = 17045.344
$ws.Cells.Item(136, 9).Value = 1882.7693
$ws.Cells.Item(136, 10).Value = 27419.736
$ws.Cells.Item(136, 11).Value = 5648.3079
$ws.Cells.Item(136, 12).Value = 82259.208
$ws.Cells.Item(136, 13).Value = -3098.3079
$ws.Cells.Item(136, 14).Value = -87359.208

$ws = $wb.Worksheets.Item("CUL")
# Row 24
$ws.Cells.Item(24, 8).Value = 791.5
$ws.Cells.Item(24, 10).Value = 791.5
$ws.Cells.Item(24, 12).Value = 2374.5
$ws.Cells.Item(24, 14).Value = -2834.5
# Row 70
$ws.Cells.Item(70, 8).Value = 4398
$ws.Cells.Item(70, 10).Value = 5012
$ws.Cells.Item(70, 12).Value = 15036
$ws.Cells.Item(70, 14).Value = -15666
# Row 73
$ws.Cells.Item(73, 8).Value = 4398
$ws.Cells.Item(73, 10).Value = 5012
$ws.Cells.Item(73, 12).Value = 15036
$ws.Cells.Item(73, 14).Value = -17220
# Row 81
$ws.Cells.Item(81, 8).Value = 5347.857
$ws.Cells.Item(81, 10).Value = 5347.857
$ws.Cells.Item(81, 12).Value = 16043.571
$ws.Cells.Item(81, 14).Value = -18289.571
# Row 84
$ws.Cells.Item(84, 8).Value = 5347.857
$ws.Cells.Item(84, 10).Value = 5347.857
$ws.Cells.Item(84, 12).Value = 48130.713
$ws.Cells.Item(84, 14).Value = -59362.713
# Row 103
$ws.Cells.Item(103, 8).Value = 1871.3077
$ws.Cells.Item(103, 9).Value = 900
$ws.Cells.Item(103, 11).Value = 2700
$ws.Cells.Item(103, 13).Value = -1821
# Row 117
$ws.Cells.Item(117, 8).Value = 1419.0667
$ws.Cells.Item(117, 9).Value = 997.8333
$ws.Cells.Item(117, 10).Value = 1699.8889
$ws.Cells.Item(117, 11).Value = 2993.4999
$ws.Cells.Item(117, 12).Value = 5099.6667
$ws.Cells.Item(117, 13).Value = 448.5001000000002
$ws.Cells.Item(117, 14).Value = -11983.6667
# Row 122
$ws.Cells.Item(122, 8).Value = 929.5
$ws.Cells.Item(122, 10).Value = 967.2222
$ws.Cells.Item(122, 12).Value = 8704.9998
$ws.Cells.Item(122, 14).Value = -13604.9998
# Row 131
$ws.Cells.Item(131, 8).Value = 752.35
$ws.Cells.Item(131, 9).Value = 0
$ws.Cells.Item(131, 10).Value = 752.35
$ws.Cells.Item(131, 11).Value = 0
$ws.Cells.Item(131, 12).Value = 2257.05
$ws.Cells.Item(131, 13).ClearContents()
$ws.Cells.Item(131, 14).Value = -12337.05

$ws = $wb.Worksheets.Item("GSM")
# Row 132
$ws.Cells.Item(132, 8).Value = 25948.375
$ws.Cells.Item(132, 9).Value = 5413.7896
$ws.Cells.Item(132, 11).Value = 16241.3688
$ws.Cells.Item(132, 13).Value = -13711.3688

$ws = $wb.Worksheets.Item("LTW")
# Row 68
$ws.Cells.Item(68, 8).Value = 2197.7144
$ws.Cells.Item(68, 9).Value = 1129.3334
$ws.Cells.Item(68, 11).Value = 1129.3334
$ws.Cells.Item(68, 13).Value = -380.3334
# Row 71
$ws.Cells.Item(71, 8).Value = 2197.7144
$ws.Cells.Item(71, 9).Value = 1129.3334
$ws.Cells.Item(71, 11).Value = 5646.666999999999
$ws.Cells.Item(71, 13).Value = -1902.666999999999
# Row 132
$ws.Cells.Item(132, 8).Value = 433062.7
$ws.Cells.Item(132, 9).Value = 503948.22
$ws.Cells.Item(132, 11).Value = 1511844.66
$ws.Cells.Item(132, 13).Value = -1509314.66

$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Cells.Item(132, 8).Value = 1345.3793
$ws.Cells.Item(132, 9).Value = 1079.4546
$ws.Cells.Item(132, 10).Value = 1507.8889
$ws.Cells.Item(132, 11).Value = 3238.3638
$ws.Cells.Item(132, 12).Value = 4523.6667
$ws.Cells.Item(132, 13).Value = -708.3638000000001
$ws.Cells.Item(132, 14).Value = -9583.6667
# Row 136
$ws.Cells.Item(136, 8).Value = 32261528
$ws.Cells.Item(136, 9).Value = 44882390
$ws.Cells.Item(136, 11).Value = 134647170
$ws.Cells.Item(136, 13).Value = -134644620
